$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.603.14'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '1.883.90'
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("D4").Value = '''1.0000'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''246.25'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '''1.0000'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '''0.2894'
$ws.Range("E8").Value = '  -1.05%  '
$ws.Range("D9").Value = '''0.06545'
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").Value = '''22.36'
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("E11").Value = '  +2.98%  '
$ws.Range("D12").Value = '''99.76'
$ws.Range("E12").Value = '  +3.19%  '
$ws.Range("D13").Value = '''0.07820'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").Value = '1.882.42'
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").Value = '''5.240'
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '''284.42'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '30.567.78'
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = '''0.000007519'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = '2.126.02'
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").Value = '''5.355'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").Value = '''0.9998'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '''6.435'
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("D25").Value = '''9.179'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").Value = '''163.64'
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = '''19.03'
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").Value = '''1.908'
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").Value = '''0.09763'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("D31").Value = '''1.503'
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("D32").Value = '''4.256'
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("D33").Value = '''4.190'
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").Value = '''0.04851'
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = '''0.7008'
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").Value = '''2.774'
$ws.Range("E37").Value = '  +2.48%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("D39").Value = '''2.867'
$ws.Range("E39").Value = '  +1.07%  '
$ws.Range("D40").Value = '''6.322'
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").Value = '''75.63'
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").Value = '''1.979'
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("D43").Value = '''0.4251'
$ws.Range("E43").Value = '  -1.05%  '
$ws.Range("D45").Value = '''0.8382'
$ws.Range("E45").Value = '  +0.59%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''9.919'
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''101.36'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("D48").Value = '''7.022'
$ws.Range("D49").Value = '''35.29'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").Value = '''0.05780'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '''0.3960'
$ws.Range("E51").Value = '  -0.45%  '
